# Generate Report for Handoff
# Updates the "Ready for handoff" rows (4-7) on the zh-cn and de-de sheets:
#  - Priority column (E) goes from "low" to "ht"
#  - Latest Handoff Datetime column (H) is bumped to reflect the newer
#    handoff-xliff-generation run.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($row in 4..7) {
    $ws_zhcn.Cells.Item($row, 5).Value = "ht"
    $ws_zhcn.Cells.Item($row, 8).Value = "2016-08-28 16:32:26"
}

$ws_dede = $wb.Worksheets.Item("de-de")
foreach ($row in 4..7) {
    $ws_dede.Cells.Item($row, 5).Value = "ht"
    $ws_dede.Cells.Item($row, 8).Value = "2016-08-28 16:32:30"
}
